# -Data de frequencia agora começa automaticamente no dia atual
# -Inicio do trabalho no sistema de procentagem
#
# Appends 22 new attendance rows (rows 14-35) for "Maria Antonieta"
# (matricula 111), mirroring the existing rows for the other students,
# and extends the used range/dimension accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$matricula = "111"
$nome = "Maria Antonieta"
$presenca = "P"

# Rows 14-29 keep the placeholder date already used elsewhere in the sheet.
for ($r = 14; $r -le 29; $r++) {
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $matricula
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $nome
    $ws.Cells.Item($r, 3).NumberFormat = "@"
    $ws.Cells.Item($r, 3).Value = "2000-01-01"
    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 4).Value = $presenca
}

# Rows 30-35: attendance now automatically starts on "today" -
# recorded here with the current date at the time of the edit.
for ($r = 30; $r -le 35; $r++) {
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $matricula
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $nome
    $ws.Cells.Item($r, 3).NumberFormat = "@"
    $ws.Cells.Item($r, 3).Value = "2024-09-21"
    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 4).Value = $presenca
}

# Last row keeps the trailing empty Justificativa/Observacoes cells,
# same as the previous last row (13) used to (present, but blank, text
# cells - a leading apostrophe forces an explicit empty-text cell rather
# than leaving it unset).
$ws.Cells.Item(35, 5).Value = "'"
$ws.Cells.Item(35, 6).Value = "'"
